# Update SQL queries and add optional supplier number parameter
#
# Data-visible effects on "Table 1":
#  - Column I ("Genuine"/"Aftermarket" supplier type) for the first 9 data
#    rows (rows 2-10) is changed from "Genuine" to "Aftermarket", adding a
#    new shared string to the workbook.
#  - The active selection is moved to K10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 10; $r++) {
    $ws.Range("I$r").Value = "Aftermarket"
}

$ws.Range("K10").Select()
